$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the error columns
$ws.Range("I1").Value = "errorDeparture"
$ws.Range("J1").Value = "errorLatitude"

# Per-station error of departure / error of latitude values (rows 2-6)
$ws.Range("I2").Value = 0.06880548829701515
$ws.Range("J2").Value = 0.1038337368845841

$ws.Range("I3").Value = 0.08878127522195502
$ws.Range("J3").Value = 0.1339790153349472

$ws.Range("I4").Value = 0.110532687651334
$ws.Range("J4").Value = 0.1668038740920093

$ws.Range("I5").Value = 0.08434221146085727
$ws.Range("J5").Value = 0.1272800645681999

$ws.Range("I6").Value = 0.1975383373688499
$ws.Range("J6").Value = 0.2981033091202576

# The tiny rounding-error total that used to live in F8 moves up to F7
# (value + formatting both copied), and F8 is then wiped out completely.
$ws.Range("F8").Copy($ws.Range("F7"))
$ws.Range("F8").Clear()

# Totals of the new error columns
$ws.Range("I8").Value = 0.5500000000000114
$ws.Range("J8").Value = 0.8299999999999983

# Precision ratio text added below the totals row
$ws.Range("B9").Value = " 1 / 1245"
